# Updates cryptos list values (price and 1h volume change) to reflect
# the latest scrape, and fixes the OKB/Stacks row order (rows 49-50 swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "93.630.93"
$ws.Range("E2").Value = "  -0.26%  "
# Row 3
$ws.Range("D3").Value = "3.431.49"
# Row 4
$ws.Range("E4").Value = "  -0.01%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.61"
$ws.Range("E5").Value = "  -0.61%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "622.35"
$ws.Range("E6").Value = "  -3.01%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.41"
$ws.Range("E7").Value = "  -1.10%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.398"
$ws.Range("E8").Value = "  -0.08%  "
# Row 9
$ws.Range("E9").Value = "  +0.02%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.984"
$ws.Range("E10").Value = "  +2.51%  "
# Row 11
$ws.Range("D11").Value = "3.433.14"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.34"
$ws.Range("E12").Value = "  +5.22%  "
# Row 13
$ws.Range("E13").Value = "  +0.83%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.33"
$ws.Range("E14").Value = "  +2.86%  "
# Row 15
$ws.Range("D15").Value = "93.490.80"
$ws.Range("E15").Value = "  -0.31%  "
# Row 16
$ws.Range("D16").Value = "4.077.05"
$ws.Range("E16").Value = "  +1.23%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000248"
$ws.Range("E17").Value = "  +0.22%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.28"
$ws.Range("E18").Value = "  +0.70%  "
# Row 19
$ws.Range("D19").Value = "3.421.94"
$ws.Range("E19").Value = "  +0.94%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.16"
$ws.Range("E20").Value = "  +5.26%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.69"
$ws.Range("E21").Value = "  +2.64%  "
# Row 22
$ws.Range("E22").Value = "  +5.36%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "504.79"
$ws.Range("E23").Value = "  +1.85%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.469"
$ws.Range("E24").Value = "  +1.07%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.67"
$ws.Range("E25").Value = "  +4.05%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000186"
$ws.Range("E26").Value = "  -1.97%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.44"
$ws.Range("E27").Value = "  +5.12%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.05"
$ws.Range("E28").Value = "  +2.39%  "
# Row 29
$ws.Range("D29").Value = "3.608.99"
$ws.Range("E29").Value = "  +0.92%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.47"
$ws.Range("E30").Value = "  +0.34%  "
# Row 31
$ws.Range("E31").Value = "  +0.01%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.140"
$ws.Range("E32").Value = "  +3.29%  "
# Row 33
$ws.Range("E33").Value = "  +1.73%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.01"
$ws.Range("E34").Value = "  +0.45%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.174"
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "29.97"
$ws.Range("E36").Value = "  +2.63%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.558"
$ws.Range("E37").Value = "  +3.24%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "560.49"
$ws.Range("E38").Value = "  +4.64%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.52"
$ws.Range("E39").Value = "  -0.91%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.40"
$ws.Range("E40").Value = "  -1.49%  "
# Row 41
$ws.Range("E41").Value = "  -0.02%  "
# Row 42
$ws.Range("E42").Value = "  +0.68%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.914"
$ws.Range("E43").Value = "  +2.14%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.74"
$ws.Range("E44").Value = "  +2.50%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.69"
$ws.Range("E45").Value = "  -1.39%  "
# Row 46
$ws.Range("E46").Value = "  +0.07%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0416"
$ws.Range("E47").Value = "  +3.77%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.53"
$ws.Range("E48").Value = "  -1.72%  "
# Row 49
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.14"
$ws.Range("E49").Value = "  -1.68%  "
# Row 50
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.66"
$ws.Range("E50").Value = "  -0.73%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.18"
$ws.Range("E51").Value = "  +2.58%  "
